$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "31.043.61"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "1.978.02"
$ws.Range("E3").Value = "  +0.29%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +1.67%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.37%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7251"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.66%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.46%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3347"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.43%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.53"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.13%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07073"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8224"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08081"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.75%  "

$ws.Range("D13").Value = "1.979.73"
$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.548"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.18%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "98.53"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +9.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "267.20"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.52%  "

$ws.Range("D18").Value = "31.040.37"
$ws.Range("E18").Value = "  +0.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.029"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.06%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000008151"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.83%  "

$ws.Range("D21").Value = "2.243.80"
$ws.Range("E21").Value = "  +1.98%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.005"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.007"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.75%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.028"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.900"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.93%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "161.78"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.85%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.50%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.338"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1322"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.67%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.583"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.372"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.64%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.589"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.11%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.382"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.17%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05263"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.18%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.270"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +5.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7724"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.37%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.785"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.10%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01992"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.98%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.870"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.71%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "82.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.75%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.709"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.21%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4601"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.48%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.079"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.41%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8510"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.005"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "104.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.08"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.608"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.52%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.569"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.05"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4266"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.69%  "
